# Auto-generated Excel COM-interop script
# Applies scheduled-runner price/profit updates to the Kujata Sheets workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9 (ALC)
$ws.Range("H9").Value = 154.75
$ws.Range("I9").Value = 175
$ws.Range("J9").Value = 134.5
$ws.Range("K9").Value = 175
$ws.Range("L9").Value = 134.5
$ws.Range("M9").Value = -6
$ws.Range("N9").Value = -472.5

# Row 33 (ALC)
$ws.Range("H33").Value = 458.8095
$ws.Range("J33").Value = 494.2
$ws.Range("L33").Value = 494.2
$ws.Range("N33").Value = -952.2

# Row 69 (ALC)
$ws.Range("H69").Value = 3006.5
$ws.Range("J69").Value = 4000
$ws.Range("L69").Value = 12000
$ws.Range("N69").Value = -13748

# Row 72 (ALC)
$ws.Range("H72").Value = 3006.5
$ws.Range("J72").Value = 4000
$ws.Range("L72").Value = 36000
$ws.Range("N72").Value = -44736

# Row 76 (ALC)
$ws.Range("H76").Value = 5800
$ws.Range("I76").Value = 6000
$ws.Range("K76").Value = 6000
$ws.Range("M76").Value = -5685

# Row 79 (ALC)
$ws.Range("H79").Value = 5800
$ws.Range("I79").Value = 6000
$ws.Range("K79").Value = 6000
$ws.Range("M79").Value = -4908

# Row 80 (ALC)
$ws.Range("H80").Value = 1331
$ws.Range("I80").Value = 2600
$ws.Range("J80").Value = 968.4286
$ws.Range("K80").Value = 7800
$ws.Range("L80").Value = 2905.2858
$ws.Range("M80").Value = -6802
$ws.Range("N80").Value = -4901.2858

# Row 83 (ALC)
$ws.Range("H83").Value = 1331
$ws.Range("I83").Value = 2600
$ws.Range("J83").Value = 968.4286
$ws.Range("K83").Value = 23400
$ws.Range("L83").Value = 8715.857399999999
$ws.Range("M83").Value = -18408
$ws.Range("N83").Value = -18699.8574

# Row 94 (ALC)
$ws.Range("H94").Value = 2500
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()

# Row 132 (ALC)
$ws.Range("H132").Value = 7941776
$ws.Range("I132").Value = 11117089
$ws.Range("J132").Value = 3493.0833
$ws.Range("K132").Value = 33351267
$ws.Range("L132").Value = 10479.2499
$ws.Range("M132").Value = -33348737
$ws.Range("N132").Value = -15539.2499

# Row 140 (ALC)
$ws.Range("H140").Value = 34425.832
$ws.Range("J140").Value = 34425.832
$ws.Range("L140").Value = 34425.832
$ws.Range("N140").Value = -44785.832

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (ARM)
$ws.Range("H32").Value = 4962.28
$ws.Range("I32").Value = 3915.0134
$ws.Range("J32").Value = 7942.9614
$ws.Range("K32").Value = 3915.0134
$ws.Range("L32").Value = 7942.9614
$ws.Range("M32").Value = -3628.0134
$ws.Range("N32").Value = -8516.9614

# Row 63 (ARM)
$ws.Range("H63").Value = 21741056
$ws.Range("I63").Value = 1803.25
$ws.Range("K63").Value = 1803.25
$ws.Range("M63").Value = -1117.25

# Row 66 (ARM)
$ws.Range("H66").Value = 21741056
$ws.Range("I66").Value = 1803.25
$ws.Range("K66").Value = 9016.25
$ws.Range("M66").Value = -5584.25

# Row 74 (ARM)
$ws.Range("H74").Value = 3069
$ws.Range("I74").Value = 2850
$ws.Range("K74").Value = 2850
$ws.Range("M74").Value = -1976

# Row 77 (ARM)
$ws.Range("H77").Value = 3069
$ws.Range("I77").Value = 2850
$ws.Range("K77").Value = 14250
$ws.Range("M77").Value = -9882

# Row 122 (ARM)
$ws.Range("H122").Value = 2212.6875
$ws.Range("I122").Value = 1881.3334
$ws.Range("J122").Value = 3206.75
$ws.Range("K122").Value = 5644.0002
$ws.Range("L122").Value = 9620.25
$ws.Range("M122").Value = -3194.0002
$ws.Range("N122").Value = -14520.25

# Row 132 (ARM)
$ws.Range("H132").Value = 2522.0232
$ws.Range("I132").Value = 2337.762
$ws.Range("K132").Value = 7013.286
$ws.Range("M132").Value = -4483.286

$ws = $wb.Worksheets.Item("CRP")
# Row 35 (CRP)
$ws.Range("H35").Value = 500
$ws.Range("I35").Value = 500
$ws.Range("K35").Value = 500
$ws.Range("M35").Value = -206

# Row 58 (CRP)
$ws.Range("H58").Value = 4010.6584
$ws.Range("I58").Value = 1186.4
$ws.Range("J58").Value = 8423.5625
$ws.Range("K58").Value = 1186.4
$ws.Range("L58").Value = 8423.5625
$ws.Range("M58").Value = -983.4000000000001
$ws.Range("N58").Value = -8829.5625

# Row 76 (CRP)
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()

# Row 79 (CRP)
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()

# Row 107 (CRP)
$ws.Range("H107").Value = 603.10345
$ws.Range("J107").Value = 1200.375
$ws.Range("L107").Value = 1200.375
$ws.Range("N107").Value = -5040.375

# Row 136 (CRP)
$ws.Range("H136").Value = 4010.6584
$ws.Range("I136").Value = 1186.4
$ws.Range("J136").Value = 8423.5625
$ws.Range("K136").Value = 3559.2
$ws.Range("L136").Value = 25270.6875
$ws.Range("M136").Value = -1009.2
$ws.Range("N136").Value = -30370.6875

$ws = $wb.Worksheets.Item("CUL")
# Row 38 (CUL)
$ws.Range("H38").Value = 100.6
$ws.Range("I38").Value = 65.75
$ws.Range("K38").Value = 197.25
$ws.Range("M38").Value = 149.75

# Row 87 (CUL)
$ws.Range("H87").Value = 2014
$ws.Range("I87").Value = 2014
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 6042
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = -4794
$ws.Range("N87").ClearContents()

# Row 90 (CUL)
$ws.Range("H90").Value = 2014
$ws.Range("I90").Value = 2014
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 18126
$ws.Range("L90").Value = 0
$ws.Range("M90").Value = -11886
$ws.Range("N90").ClearContents()

# Row 98 (CUL)
$ws.Range("H98").Value = 266.66666
$ws.Range("I98").Value = 181
$ws.Range("J98").Value = 283.8
$ws.Range("K98").Value = 543
$ws.Range("L98").Value = 851.4000000000001
$ws.Range("M98").Value = 955
$ws.Range("N98").Value = -3847.4

# Row 107 (CUL)
$ws.Range("H107").Value = 4359.76
$ws.Range("I107").Value = 400
$ws.Range("J107").Value = 4899.727
$ws.Range("K107").Value = 1200
$ws.Range("L107").Value = 14699.181
$ws.Range("M107").Value = 720
$ws.Range("N107").Value = -18539.181

# Row 113 (CUL)
$ws.Range("H113").Value = 690.4194
$ws.Range("I113").Value = 586.6
$ws.Range("J113").Value = 739.8570999999999
$ws.Range("K113").Value = 1759.8
$ws.Range("L113").Value = 2219.5713
$ws.Range("M113").Value = 410.1999999999998
$ws.Range("N113").Value = -6559.5713

# Row 131 (CUL)
$ws.Range("H131").Value = 16950598
$ws.Range("J131").Value = 1666.8959
$ws.Range("L131").Value = 5000.6877
$ws.Range("N131").Value = -15080.6877

# Row 132 (CUL)
$ws.Range("H132").Value = 1074.25
$ws.Range("I132").Value = 1042
$ws.Range("J132").Value = 1300
$ws.Range("K132").Value = 9378
$ws.Range("L132").Value = 11700
$ws.Range("M132").Value = -6848
$ws.Range("N132").Value = -16760

# Row 137 (CUL)
$ws.Range("H137").Value = 34097844
$ws.Range("I137").Value = 53573212
$ws.Range("J137").Value = 15945.75
$ws.Range("K137").Value = 160719636
$ws.Range("L137").Value = 47837.25
$ws.Range("M137").Value = -160714536
$ws.Range("N137").Value = -58037.25

$ws = $wb.Worksheets.Item("GSM")
# Row 93 (GSM)
$ws.Range("H93").Value = 29960
$ws.Range("J93").Value = 29960
$ws.Range("L93").Value = 29960
$ws.Range("N93").Value = -33704

# Row 119 (GSM)
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()

# Row 132 (GSM)
$ws.Range("H132").Value = 3113.889
$ws.Range("I132").Value = 2818.15
$ws.Range("J132").Value = 3958.8572
$ws.Range("K132").Value = 8454.450000000001
$ws.Range("L132").Value = 11876.5716
$ws.Range("M132").Value = -5924.450000000001
$ws.Range("N132").Value = -16936.5716

$ws = $wb.Worksheets.Item("LTW")
# Row 13 (LTW)
$ws.Range("H13").Value = 2333.3333
$ws.Range("I13").Value = 2333.3333
$ws.Range("K13").Value = 2333.3333
$ws.Range("M13").Value = -2193.3333

# Row 132 (LTW)
$ws.Range("H132").Value = 2486.2666
$ws.Range("I132").Value = 2450.1875
$ws.Range("J132").Value = 2527.5
$ws.Range("K132").Value = 7350.5625
$ws.Range("L132").Value = 7582.5
$ws.Range("M132").Value = -4820.5625
$ws.Range("N132").Value = -12642.5

$ws = $wb.Worksheets.Item("WVR")
# Row 113 (WVR)
$ws.Range("H113").Value = 401.70587
$ws.Range("I113").Value = 331.7
$ws.Range("J113").Value = 501.7143
$ws.Range("K113").Value = 995.0999999999999
$ws.Range("L113").Value = 1505.1429
$ws.Range("M113").Value = 1174.9
$ws.Range("N113").Value = -5845.1429

# Row 132 (WVR)
$ws.Range("H132").Value = 2966.682
$ws.Range("I132").Value = 3776.4
$ws.Range("J132").Value = 2291.9167
$ws.Range("K132").Value = 11329.2
$ws.Range("L132").Value = 6875.750100000001
$ws.Range("M132").Value = -8799.200000000001
$ws.Range("N132").Value = -11935.7501
